$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-98 down to 16-99
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new record's data
$ws.Cells.Item(15, 1).Value = 11
$ws.Cells.Item(15, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(15, 3).Value = "Bíobío"
$ws.Cells.Item(15, 4).Value = 44749
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 6).Value = 100112001
$ws.Cells.Item(15, 7).Value = "Berenjena"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 170
$ws.Cells.Item(15, 11).Value = 13000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 13941
$ws.Cells.Item(15, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 232
$ws.Cells.Item(15, 17).Value = 60
$ws.Cells.Item(15, 18).Value = "Hortaliza"
